$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D stays text-formatted so numeric-looking price strings
# (e.g. "1.00", "0.180") are not silently coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "66.261.49"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "3.352.23"
$ws.Range("E3").Value = "  +2.72%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "189.49"
$ws.Range("E5").Value = "  +4.80%  "
$ws.Range("D6").Value = "558.78"
$ws.Range("E6").Value = "  +0.60%  "
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("D8").Value = "3.346.10"
$ws.Range("E8").Value = "  +2.77%  "
$ws.Range("D9").Value = "0.584"
$ws.Range("E9").Value = "  -0.69%  "
$ws.Range("E10").Value = "  -2.87%  "
$ws.Range("D11").Value = "0.583"
$ws.Range("E11").Value = "  -0.38%  "
$ws.Range("D12").Value = "46.39"
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("E13").Value = "  +1.61%  "
$ws.Range("D14").Value = "3.887.92"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").Value = "8.55"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").Value = "592.29"
$ws.Range("E16").Value = "  -6.74%  "
$ws.Range("D17").Value = "66.263.77"
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "3.344.32"
$ws.Range("E18").Value = "  +2.75%  "
$ws.Range("D19").Value = "17.94"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("E20").Value = "  +1.16%  "
$ws.Range("D21").Value = "11.03"
$ws.Range("E21").Value = "  -2.75%  "
$ws.Range("D22").Value = "0.899"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "18.24"
$ws.Range("E23").Value = "  +2.83%  "
$ws.Range("D25").Value = "99.39"
$ws.Range("E25").Value = "  -6.17%  "
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "6.03"
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").Value = "2.72"
$ws.Range("E28").Value = "  +2.21%  "
$ws.Range("D29").Value = "9.46"
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("D32").Value = "6.69"
$ws.Range("E32").Value = "  +6.10%  "
$ws.Range("B33").Value = "Bittensor"
$ws.Range("C33").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D33").Value = "585.18"
$ws.Range("E33").Value = "  +6.30%  "
$ws.Range("B34").Value = "dogwifhat"
$ws.Range("C34").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D34").Value = "3.81"
$ws.Range("E34").Value = "  -7.03%  "
$ws.Range("D35").Value = "10.94"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").Value = "3.777.49"
$ws.Range("E37").Value = "  +4.18%  "
$ws.Range("D38").Value = "0.999"
$ws.Range("D39").Value = "56.31"
$ws.Range("E39").Value = "  -1.28%  "
$ws.Range("D40").Value = "34.23"
$ws.Range("E40").Value = "  +7.26%  "
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0699"
$ws.Range("E41").Value = "  -2.48%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.126"
$ws.Range("E42").Value = "  -4.39%  "
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("D44").Value = "3.16"
$ws.Range("E44").Value = "  -7.28%  "
$ws.Range("D45").Value = "3.41"
$ws.Range("E45").Value = "  +2.68%  "
$ws.Range("D46").Value = "0.339"
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("D47").Value = "0.0414"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "3.06"
$ws.Range("E48").Value = "  -16.75%  "
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").Value = "2.55"
$ws.Range("E50").Value = "  -2.29%  "
$ws.Range("E51").Value = "  +0.03%  "

# Restore default cell style (the text NumberFormat above is transient and
# must not remain applied once values are written, to match original styling).
$ws.Range("D2:D51").Style = "Normal"
